$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New NATMI LR-pair table for Icam5-Itgal across ECs/FAPs/M2 sending clusters
# and FAPs/M2/sCs target clusters (rows 2-10), replacing the old 2-row table.
$row2 = @("ECs","Icam5","Itgal","FAPs",2,0.6666666666666666,0.485621,1.456863,0.1357088875882878,0.1357088875882878,3,1,0.9976189999999999,2.992857,0.06139866523632875,0.06139866523632876,0.484464736399,4.360182627591,0.008332344558627851,0.008332344558627851)
$row3 = @("ECs","Icam5","Itgal","M2",2,0.6666666666666666,0.485621,1.456863,0.1357088875882878,0.1357088875882878,3,1,14.969586,44.908758,0.9213062296732858,0.9213062296732859,7.269545322906001,65.42590790615399,0.1250294435571212,0.1250294435571212)
$row4 = @("ECs","Icam5","Itgal","sCs",2,0.6666666666666666,0.485621,1.456863,0.1357088875882878,0.1357088875882878,2,0.6666666666666666,0.2810146666666667,0.843044,0.01729510509038539,0.01729510509038539,0.1364666234413333,1.228199610972,0.002347099472538735,0.002347099472538734)
$row5 = @("FAPs","Icam5","Itgal","FAPs",3,1,2.580934333333333,7.742803,0.7212532557592974,0.7212532557592974,3,1,0.9976189999999999,2.992857,0.06139866523632875,0.06139866523632876,2.574789128685667,23.173102158171,0.04428398720097731,0.04428398720097732)
$row6 = @("FAPs","Icam5","Itgal","M2",3,1,2.580934333333333,7.742803,0.7212532557592974,0.7212532557592974,3,1,14.969586,44.908758,0.9213062296732858,0.9213062296732859,38.635518463186,347.719666168674,0.6644951177031804,0.6644951177031805)
$row7 = @("FAPs","Icam5","Itgal","sCs",3,1,2.580934333333333,7.742803,0.7212532557592974,0.7212532557592974,2,0.6666666666666666,0.2810146666666667,0.843044,0.01729510509038539,0.01729510509038539,0.7252804013702223,6.527523612332001,0.01247415085513966,0.01247415085513966)
$row8 = @("M2","Icam5","Itgal","FAPs",2,0.6666666666666666,0.5118469999999999,1.535541,0.1430378566524148,0.1430378566524148,3,1,0.9976189999999999,2.992857,0.06139866523632875,0.06139866523632876,0.5106282922929999,4.595654630636999,0.008782333476723594,0.008782333476723596)
$row9 = @("M2","Icam5","Itgal","M2",2,0.6666666666666666,0.5118469999999999,1.535541,0.1430378566524148,0.1430378566524148,3,1,14.969586,44.908758,0.9213062296732858,0.9213062296732859,7.662137685341999,68.95923916807799,0.1317816684129842,0.1317816684129842)
$row10 = @("M2","Icam5","Itgal","sCs",2,0.6666666666666666,0.5118469999999999,1.535541,0.1430378566524148,0.1430378566524148,2,0.6666666666666666,0.2810146666666667,0.843044,0.01729510509038539,0.01729510509038539,0.1438365140893333,1.294528626804,0.002473854762706995,0.002473854762706995)

$rows = @($row2, $row3, $row4, $row5, $row6, $row7, $row8, $row9, $row10)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $vals = $rows[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value = $vals[$j]
    }
}
